$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 5000
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("N3").Value = -5228
$ws.Range("H64").Value = 3519.1428
$ws.Range("J64").Value = 3509.6956
$ws.Range("L64").Value = 3509.6956
$ws.Range("N64").Value = -4005.6956
$ws.Range("H67").Value = 3519.1428
$ws.Range("J67").Value = 3509.6956
$ws.Range("L67").Value = 3509.6956
$ws.Range("N67").Value = -5225.6956
$ws.Range("H102").Value = 5000
$ws.Range("J102").Value = 5000
$ws.Range("L102").Value = 5000
$ws.Range("N102").Value = -11490
$ws.Range("H137").Value = 1549.7878
$ws.Range("I137").Value = 1318.2174
$ws.Range("J137").Value = 2082.4
$ws.Range("K137").Value = 3954.6522
$ws.Range("L137").Value = 6247.200000000001
$ws.Range("M137").Value = -1404.6522
$ws.Range("N137").Value = -11347.2
$ws.Range("H138").Value = 574004.5
$ws.Range("I138").Value = 1785
$ws.Range("J138").Value = 650869.8
$ws.Range("K138").Value = 5355
$ws.Range("L138").Value = 1952609.4
$ws.Range("M138").Value = -215
$ws.Range("N138").Value = -1962889.4
$ws.Range("H141").Value = 1673.3125
$ws.Range("I141").Value = 1673.3125
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5019.9375
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 160.0625
$ws.Range("N141").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2318.8333
$ws.Range("I61").Value = 2133
$ws.Range("J61").Value = 2504.6667
$ws.Range("K61").Value = 2133
$ws.Range("L61").Value = 2504.6667
$ws.Range("M61").Value = -1921
$ws.Range("N61").Value = -2928.6667
$ws.Range("H74").Value = 1755.875
$ws.Range("I74").Value = 1674.5
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 1674.5
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -800.5
$ws.Range("N74").Value = -3748
$ws.Range("H77").Value = 1755.875
$ws.Range("I77").Value = 1674.5
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 8372.5
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -4004.5
$ws.Range("N77").Value = -18736
$ws.Range("H88").Value = 2643.0833
$ws.Range("I88").Value = 1003
$ws.Range("J88").Value = 2971.1
$ws.Range("K88").Value = 1003
$ws.Range("L88").Value = 2971.1
$ws.Range("M88").Value = -597
$ws.Range("N88").Value = -3783.1
$ws.Range("H91").Value = 2643.0833
$ws.Range("I91").Value = 1003
$ws.Range("J91").Value = 2971.1
$ws.Range("K91").Value = 1003
$ws.Range("L91").Value = 2971.1
$ws.Range("M91").Value = 401
$ws.Range("N91").Value = -5779.1
$ws.Range("H132").Value = 2908.0952
$ws.Range("I132").Value = 2789.9285
$ws.Range("K132").Value = 8369.7855
$ws.Range("M132").Value = -5839.7855
$ws.Range("H135").Value = 16947.334
$ws.Range("J135").Value = 16947.334
$ws.Range("L135").Value = 16947.334
$ws.Range("N135").Value = -27087.334
$ws.Range("H136").Value = 2318.8333
$ws.Range("I136").Value = 2133
$ws.Range("J136").Value = 2504.6667
$ws.Range("K136").Value = 6399
$ws.Range("L136").Value = 7514.000100000001
$ws.Range("M136").Value = -3849
$ws.Range("N136").Value = -12614.0001
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 55556690
$ws.Range("I99").Value = 55556690
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 55556690
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -55555192
$ws.Range("H134").Value = 5287.2085
$ws.Range("I134").Value = 1094.4445
$ws.Range("J134").Value = 17865.5
$ws.Range("K134").Value = 3283.3335
$ws.Range("L134").Value = 53596.5
$ws.Range("M134").Value = -748.3335000000002
$ws.Range("N134").Value = -58666.5
$ws.Range("N99").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1665.3448
$ws.Range("I31").Value = 1341.8462
$ws.Range("J31").Value = 1928.1875
$ws.Range("K31").Value = 1341.8462
$ws.Range("L31").Value = 1928.1875
$ws.Range("M31").Value = -1046.8462
$ws.Range("N31").Value = -2518.1875
$ws.Range("H34").Value = 1665.3448
$ws.Range("I34").Value = 1341.8462
$ws.Range("J34").Value = 1928.1875
$ws.Range("K34").Value = 1341.8462
$ws.Range("L34").Value = 1928.1875
$ws.Range("M34").Value = -1139.8462
$ws.Range("N34").Value = -2332.1875
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("H58").Value = 797.3077
$ws.Range("I58").Value = 795.9167
$ws.Range("J58").Value = 814
$ws.Range("K58").Value = 795.9167
$ws.Range("L58").Value = 814
$ws.Range("M58").Value = -592.9167
$ws.Range("N58").Value = -1220
$ws.Range("H122").Value = 796.6667
$ws.Range("I122").Value = 796.6667
$ws.Range("K122").Value = 2390.0001
$ws.Range("M122").Value = 59.9998999999998
$ws.Range("H132").Value = 13970
$ws.Range("I132").Value = 13970
$ws.Range("K132").Value = 41910
$ws.Range("M132").Value = -39380
$ws.Range("H136").Value = 797.3077
$ws.Range("I136").Value = 795.9167
$ws.Range("J136").Value = 814
$ws.Range("K136").Value = 2387.7501
$ws.Range("L136").Value = 2442
$ws.Range("M136").Value = 162.2498999999998
$ws.Range("N136").Value = -7542
$ws.Range("N45").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1159.5
$ws.Range("J5").Value = 499.9375
$ws.Range("L5").Value = 1499.8125
$ws.Range("N5").Value = -1723.8125
$ws.Range("H68").Value = 2335.3845
$ws.Range("J68").Value = 2335.3845
$ws.Range("L68").Value = 7006.1535
$ws.Range("N68").Value = -8628.1535
$ws.Range("H71").Value = 2335.3845
$ws.Range("J71").Value = 2335.3845
$ws.Range("L71").Value = 21018.4605
$ws.Range("N71").Value = -29130.4605
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H135").Value = 1159.5
$ws.Range("J135").Value = 499.9375
$ws.Range("L135").Value = 4499.4375
$ws.Range("N135").Value = -9569.4375
$ws.Range("H141").Value = 2031.9
$ws.Range("I141").Value = 2031.9
$ws.Range("K141").Value = 6095.700000000001
$ws.Range("M141").Value = -915.7000000000007
$ws.Range("N125").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 3051
$ws.Range("J110").Value = 3051
$ws.Range("L110").Value = 3051
$ws.Range("N110").Value = -11231
$ws.Range("H126").Value = 2467.2964
$ws.Range("I126").Value = 1324.9286
$ws.Range("K126").Value = 3974.7858
$ws.Range("M126").Value = -1504.7858
$ws.Range("H132").Value = 3135.476
$ws.Range("I132").Value = 2402.8235
$ws.Range("K132").Value = 7208.470499999999
$ws.Range("M132").Value = -4678.470499999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2405
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 2405
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 2405
$ws.Range("N61").Value = -2809
$ws.Range("H93").Value = 1067
$ws.Range("I93").Value = 466
$ws.Range("J93").Value = 1668
$ws.Range("K93").Value = 466
$ws.Range("L93").Value = 1668
$ws.Range("M93").Value = 782
$ws.Range("N93").Value = -4164
$ws.Range("H113").Value = 2405
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2405
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2405
$ws.Range("N113").Value = -6745
$ws.Range("M61").ClearContents()
$ws.Range("M113").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 38470100
$ws.Range("J62").Value = 14625
$ws.Range("L62").Value = 14625
$ws.Range("N62").Value = -15873
$ws.Range("H65").Value = 38470100
$ws.Range("J65").Value = 14625
$ws.Range("L65").Value = 73125
$ws.Range("N65").Value = -79365
$ws.Range("H113").Value = 1200
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3600
$ws.Range("N113").Value = -7940
$ws.Range("H122").Value = 32511252
$ws.Range("I122").Value = 52014800
$ws.Range("J122").Value = 5335
$ws.Range("K122").Value = 156044400
$ws.Range("L122").Value = 16005
$ws.Range("M122").Value = -156041950
$ws.Range("N122").Value = -20905
$ws.Range("M113").ClearContents()
